$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}

for ($r = 34; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
